# Test Data changes - Companies module - 4th Dec 2023
$wb = $excel.ActiveWorkbook

# Update the Users sheet: change "Drew Koecher" to "Ayati Arvind"
$ws = $wb.Worksheets.Item("Users")
$ws.Range("A2").Value = "Ayati Arvind"

# Update the selected cell on the Users sheet to A2
$ws.Activate()
$ws.Range("A2").Select()
